# Fake-Team-Data.xlsx update:
#  - add 5 new athletes (rows 19-23) to the "Athletes" sheet
#  - add a "Student Classification" column (D) to the Athletes sheet, driven by
#    a dropdown (data validation) list
#  - add a new "Data Validation Variables" sheet holding the dropdown options
#  - bold the header rows on all three sheets

$wb = $excel.ActiveWorkbook

$coaches  = $wb.Worksheets.Item("Coaches")
$athletes = $wb.Worksheets.Item("Athletes")

# ---------------------------------------------------------------------------
# 1. New athlete rows (first name / last name / age)
# ---------------------------------------------------------------------------
$athletes.Range("A19").Value = "Theo"
$athletes.Range("B19").Value = "McIntoss"
$athletes.Range("C19").Value = 18

$athletes.Range("A20").Value = "Alexa"
$athletes.Range("B20").Value = "Amanda"
$athletes.Range("C20").Value = 17

$athletes.Range("A21").Value = "Bud"
$athletes.Range("B21").Value = "Flambeau"
$athletes.Range("C21").Value = 18

$athletes.Range("A22").Value = "Burt"
$athletes.Range("B23").Value = "Lancaster"
$athletes.Range("B22").Value = "Bogart"
$athletes.Range("A23").Value = "Humphrey"
$athletes.Range("C22").Value = 19
$athletes.Range("C23").Value = 19

# ---------------------------------------------------------------------------
# 2. New "Data Validation Variables" sheet with the dropdown source list
# ---------------------------------------------------------------------------
$dv = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$dv.Name = "Data Validation Variables"

$dv.Range("A2").Value = "Student Classification"
$dv.Range("A3").Value = "Freshman"
$dv.Range("A4").Value = "Sophomore"
$dv.Range("A5").Value = "Junior"
$dv.Range("A6").Value = "Senior"
$dv.Range("A7").Value = "Other"
$dv.Range("A8").Value = "N/A"

$dv.Range("A1").Value = "NOTE: Don't modify the data below without consulting Soham or withour understanding what it does! It is used for the simple dropdowns in the other sheets."

$dv.Range("A1:A8").HorizontalAlignment = -4131
$dv.Range("A1:A8").VerticalAlignment = -4108

$dv.Range("A1").Font.Color = 255

# ---------------------------------------------------------------------------
# 3. "Student Classification" column on the Athletes sheet
# ---------------------------------------------------------------------------
$athletes.Range("D1").Value = "Student Classification"

$athletes.Range("D2").Value = "Senior"
$athletes.Range("D3").Value = "Junior"
$athletes.Range("D4").Value = "Senior"
$athletes.Range("D5").Value = "N/A"
$athletes.Range("D6").Value = "Senior"
$athletes.Range("D7").Value = "Junior"
$athletes.Range("D8").Value = "Sophomore"
$athletes.Range("D9").Value = "Junior"
$athletes.Range("D10").Value = "Senior"
$athletes.Range("D11").Value = "Junior"
$athletes.Range("D12").Value = "Junior"
$athletes.Range("D13").Value = "Junior"
$athletes.Range("D14").Value = "Sophomore"
$athletes.Range("D15").Value = "Senior"
$athletes.Range("D16").Value = "Senior"
$athletes.Range("D17").Value = "Junior"
$athletes.Range("D18").Value = "Senior"
$athletes.Range("D19").Value = "Senior"
$athletes.Range("D20").Value = "Sophomore"
$athletes.Range("D21").Value = "Senior"
$athletes.Range("D22").Value = "Senior"
$athletes.Range("D23").Value = "Senior"

$athletes.Range("D2:D1048576").Validation.Add(3, 1, 1, "='Data Validation Variables'!`$A`$3:`$A`$8")

# ---------------------------------------------------------------------------
# 4. Bold the header rows
# ---------------------------------------------------------------------------
$coaches.Range("A1:C1").Font.Bold = $true
$athletes.Range("A1:D1").Font.Bold = $true
$dv.Range("A2").Font.Bold = $true

# ---------------------------------------------------------------------------
# 5. Selections to mirror the saved workbook state
# ---------------------------------------------------------------------------
$athletes.Range("E1").Select()
$dv.Range("A2").Select()
$coaches.Activate()
$athletes.Activate()

Write-Host "done"
